$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.712.52"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "3.585.38"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'608.74"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "'145.66"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").Value = "'0.137"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "4.190.95"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "'0.0000209"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'30.06"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "3.574.23"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "66.719.60"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'11.43"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'6.23"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'432.91"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'0.621"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").Value = "'79.20"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "3.728.52"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "'9.33"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").Value = "'8.08"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "3.579.33"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").Value = "'25.47"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'7.85"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "'5.64"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'173.77"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "'0.894"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = "  +5.88%  "
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").Value = "'25.17"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "'7.22"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "'23.60"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "'0.944"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  -1.36%  "
